# Update the pantheon stats for a handful of players (rows 3, 4, 8, 10)
# to reflect newly played games. Columns: D NBGAMES, E SOLOKILLS,
# F DUREE_GAME, G WARDS_SCORE, H WARDS_POSEES, I WARDS_DETRUITES,
# J WARDS_PINKS, K CS, L KILLS, M DEATHS, N ASSISTS,
# O KILLS_MOYENNE, P DEATHS_MOYENNE, Q ASSISTS_MOYENNE,
# R WARDS_MOYENNE, S DUREE_MOYENNE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - kulbutoké
$ws.Range("D3").Value = 32
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 15.48733333333334
$ws.Range("G3").Value = 723
$ws.Range("H3").Value = 357
$ws.Range("I3").Value = 65
$ws.Range("J3").Value = 77
$ws.Range("K3").Value = 3884
$ws.Range("L3").Value = 96
$ws.Range("M3").Value = 147
$ws.Range("N3").Value = 250
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 4.59375
$ws.Range("Q3").Value = 7.8125
$ws.Range("R3").Value = 22.59
$ws.Range("S3").Value = 29.04

# Row 4 - tomlora
$ws.Range("D4").Value = 152
$ws.Range("E4").Value = 117
$ws.Range("F4").Value = 75.01533333333337
$ws.Range("G4").Value = 3058
$ws.Range("H4").Value = 1303
$ws.Range("I4").Value = 438
$ws.Range("J4").Value = 408
$ws.Range("K4").Value = 30437
$ws.Range("L4").Value = 1205
$ws.Range("M4").Value = 913
$ws.Range("N4").Value = 1151
$ws.Range("O4").Value = 7.927631578947368
$ws.Range("P4").Value = 6.006578947368421
$ws.Range("Q4").Value = 7.572368421052632
$ws.Range("R4").Value = 20.12
$ws.Range("S4").Value = 29.61

# Row 8 - nukethestars
$ws.Range("D8").Value = 33
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 16.0685
$ws.Range("G8").Value = 2375
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 361
$ws.Range("J8").Value = 227
$ws.Range("K8").Value = 1214
$ws.Range("L8").Value = 97
$ws.Range("M8").Value = 239
$ws.Range("N8").Value = 457
$ws.Range("O8").Value = 2.939393939393939
$ws.Range("P8").Value = 7.242424242424242
$ws.Range("Q8").Value = 13.84848484848485
$ws.Range("R8").Value = 71.97
$ws.Range("S8").Value = 29.22

# Row 10 - namiyeon
$ws.Range("D10").Value = 76
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 37.59033333333333
$ws.Range("G10").Value = 5911
$ws.Range("H10").Value = 3386
$ws.Range("I10").Value = 739
$ws.Range("J10").Value = 801
$ws.Range("K10").Value = 1463
$ws.Range("L10").Value = 192
$ws.Range("M10").Value = 223
$ws.Range("N10").Value = 1221
$ws.Range("O10").Value = 2.526315789473684
$ws.Range("P10").Value = 2.934210526315789
$ws.Range("Q10").Value = 16.06578947368421
$ws.Range("R10").Value = 77.78
$ws.Range("S10").Value = 29.68
